# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages data refresh at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 65
$wsExpo.Range("F3").Value = 796
$wsExpo.Range("F6").Value = 92
$wsExpo.Range("F7").Value = 319
$wsExpo.Range("F8").Value = 4033
$wsExpo.Range("F10").Value = 4742
$wsExpo.Range("F11").Value = 526
$wsExpo.Range("F12").Value = 1193

# Sheet "全部类型" (sheet4)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 65
$wsAll.Range("F3").Value = 796
$wsAll.Range("F6").Value = 92
$wsAll.Range("F8").Value = 319
$wsAll.Range("F9").Value = 4033
$wsAll.Range("F11").Value = 4742
$wsAll.Range("F12").Value = 526
$wsAll.Range("F13").Value = 1193
